# 236-GFG-Find LCA in a Binary Tree
# Append a new question row (row 21) to the Binary Search Tree question list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: Question No / GFG-LC / Question
$ws.Range("A21").Value = 236
$ws.Range("B21").Value = "LC/GFG"
$ws.Range("C21").Value = "Find LCA in a Binary tree"

# Match the formatting used by the other "LC/GFG" rows (style index 1:
# left/top aligned, wrapped text) for the middle column.
$ws.Range("B21").HorizontalAlignment = -4131
$ws.Range("B21").VerticalAlignment = -4160
$ws.Range("B21").WrapText = $true

# Update the active selection to the new last row, as in the saved file.
$ws.Range("A21").Select()
